$d = $word.ActiveDocument

# The commit "done with use reducer" marks the last four bullet/line items of
# the useRef section's last task line plus its bonus line, and the first
# three bullet items of the useReducer ("Todo List") section as complete by
# turning their text red (FF0000), matching the highlighting convention
# already used elsewhere in the document (see the useState/useEffect/
# useContext sections above them).
#
#  - "Store the number of times the user clicked a button, but do not
#     trigger re-renders"
#  - "Bonus: Show the number of clicks without making the component
#     re-render"
#  - "Build a todo list using useReducer"
#  - "The reducer should handle adding, removing, and marking tasks as done"
#  - "Display the tasks in a list"
#
# (The final useReducer bullet, the "Clear Completed" bonus line, is left
# untouched.)

$targets = @(
    "Store the number of times the user clicked a button, but do not trigger re-renders",
    "Bonus: Show the number of clicks without making the component re-render",
    "Build a todo list using useReducer",
    "The reducer should handle adding, removing, and marking tasks as done",
    "Display the tasks in a list"
)

$red = 255

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text.Trim()
    foreach ($t in $targets) {
        if ($text -eq $t -or $text -eq ("👉 " + $t)) {
            $p.Range.Font.Color = $red
            break
        }
    }
}
